$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Report date at top of sheet (A1): 1/26/2021 -> 10/5/2021
$ws.Range("A1").Value = 44474

# "As of" date under header (C6) switches from a date value to a text label
$ws.Range("C6").Value = "August, 2021"

# Santa Barbara County / California / United States summary rows
$ws.Range("C9").Value = 802120
$ws.Range("D9").Value = 0.035593524506456
$ws.Range("E9").Value = 0.218049953836435
$ws.Range("F9").Value = 1.05058697021730008

$ws.Range("C10").Value = 700828
$ws.Range("D10").Value = 0.0257976033477605
$ws.Range("E10").Value = 0.21214238644834699
$ws.Range("F10").Value = 1.33392718129683008

$ws.Range("C11").Value = 298933
$ws.Range("D11").Value = 0.0198279891239455
$ws.Range("E11").Value = 0.16732856144297001
$ws.Range("F11").Value = 1.38128243160195008

# Cities table
$ws.Range("C14").Value = 1151915
$ws.Range("D14").Value = 0.0421231569252927
$ws.Range("E14").Value = 0.22977710780346899
$ws.Range("F14").Value = 1.16815705570249007

$ws.Range("C15").Value = 1058989
$ws.Range("D15").Value = 0.039269893049283
$ws.Range("E15").Value = 0.23247889110662401
$ws.Range("F15").Value = 1.09234934912219006

$ws.Range("C16").Value = 466456
$ws.Range("D16").Value = 0.0369698505644405
$ws.Range("E16").Value = 0.24061076254890401
$ws.Range("F16").Value = 0.93293812226742301

$ws.Range("C17").Value = 4126752
$ws.Range("D17").Value = 0.034463465477877
$ws.Range("E17").Value = 0.239796152885936
$ws.Range("F17").Value = 1.07541658935641005

$ws.Range("C18").Value = 1468019
$ws.Range("D18").Value = 0.0387987735539974
$ws.Range("E18").Value = 0.24626274580049701
$ws.Range("F18").Value = 1.15567113866772009

$ws.Range("C19").Value = 512099
$ws.Range("D19").Value = 0.0348946713227924
$ws.Range("E19").Value = 0.22210581581270999
$ws.Range("F19").Value = 1.01685831667027005

$ws.Range("C20").Value = 1042019
$ws.Range("D20").Value = 0.0318255949475135
$ws.Range("E20").Value = 0.24738466989799601
$ws.Range("F20").Value = 1.0434307618207399

# Move the active selection to match the saved view state
$ws.Range("I21").Select()
